{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies the 4 textual/formatting changes described by the diff:\n//  1. \"Monsieur le Directeur de la soci\u00e9t\u00e9\u2026\u2026\u2026\u2026\u2026\u2026\u2026.\"  -> \"...soci\u00e9t\u00e9 <societe>\"\n//  2. \"Objet    : ... n\u00b0 : \u2026\u2026\u2026\u2026.  \" -> \"Objet    : ... n\u00b0 : <Num_PR> \"\n//  3. \"Suite \u00e0 votre demande ...\" paragraph -> filled with <Num_PR>/Marrakech/<date>\n//  4. Closing \"\u2026\u2026\u2026\u2026\u2026\u2026, le\u2026\u2026\u2026/\u2026\u2026\u2026/\u2026\u2026\u2026\" line -> \"\u2026Marrakech, le <date>\" + centered\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// ---- 1. \"Monsieur le Directeur de la soci\u00e9t\u00e9\u2026\" paragraph ----------------\n// Find the short (salutation block) paragraph that starts with this text and\n// still contains the placeholder ellipsis dots, so we don't also match the\n// \"Monsieur le Directeur,\" greeting later in the letter.\nlet directeurPara = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(\"Monsieur le Directeur de la soci\u00e9t\u00e9\") !== -1) {\n    directeurPara = p;\n    break;\n  }\n}\nif (directeurPara) {\n  directeurPara.insertText(\n    \"Monsieur le Directeur de la soci\u00e9t\u00e9 <societe>\",\n    Word.InsertLocation.replace\n  );\n}\n\n// ---- 2. \"Objet\" paragraph -------------------------------------------------\n// \"Objet \" itself stays bold+underlined and untouched; only the remainder of\n// the paragraph (after that label) is replaced.\nconst objetResults = context.document.body.search(\"Objet \", { matchCase: false });\nobjetResults.load(\"items\");\nawait context.sync();\n\nif (objetResults.items.length > 0) {\n  const objetMatch = objetResults.items[0];\n  const objetPara = objetMatch.paragraphs.getFirst();\n  const afterObjet = objetMatch.getRange(\"After\");\n  const objetTail = afterObjet.expandTo(objetPara.getRange(\"End\"));\n  await context.sync();\n  objetTail.insertText(\n    \"   : D\u00e9cision d\\u2019Institution du permis de recherche n\u00b0 : <Num_PR> \",\n    Word.InsertLocation.replace\n  );\n}\n\n// ---- 3. \"Suite \u00e0 votre demande ...\" paragraph -----------------------------\nlet suitePara = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(\"Suite \\u00e0 votre demande\") !== -1) {\n    suitePara = p;\n    break;\n  }\n}\nif (suitePara) {\n  suitePara.insertText(\n    \"                Suite \\u00e0 votre demande d\\u2019Institution du permis de recherche n\u00b0 <Num_PR> \" +\n      \"d\\u00e9pos\\u00e9e \\u00e0 La Direction R\\u00e9gionale du D\\u00e9partement de l\\u2019Energie et des Mines de Marrakech. \" +\n      \"Le <date>, j\\u2019ai l\\u2019honneur de vous faire parvenir ci-joint la d\\u00e9cision cit\\u00e9e en objet. \",\n    Word.InsertLocation.replace\n  );\n}\n\n// ---- 4. Closing date/place line -------------------------------------------\nlet closingPara = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(\"\u2026\u2026\u2026\u2026\u2026\u2026, le\") !== -1) {\n    closingPara = p;\n    break;\n  }\n}\nif (closingPara) {\n  closingPara.insertText(\n    \"                                                                   Marrakech, le <date>\",\n    Word.InsertLocation.replace\n  );\n  closingPara.alignment = Word.Alignment.centered;\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop script.\n# Applies the 4 textual/formatting changes described by the diff:\n#  1. \"Monsieur le Directeur de la soci\u00e9t\u00e9\u2026\u2026\u2026\u2026\u2026\u2026\u2026.\"  -> \"...soci\u00e9t\u00e9 <societe>\"\n#  2. \"Objet    : ... n\u00b0 : \u2026\u2026\u2026\u2026.  \" -> \"Objet    : ... n\u00b0 : <Num_PR> \"\n#  3. \"Suite \u00e0 votre demande ...\" paragraph -> filled with <Num_PR>/Marrakech/<date>\n#  4. Closing \"\u2026\u2026\u2026\u2026\u2026\u2026, le\u2026\u2026\u2026/\u2026\u2026\u2026/\u2026\u2026\u2026\" line -> \"\u2026Marrakech, le <date>\" + centered\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($oldText, $newText) {\n    $find = $d.Content.Find\n    $find.Text = $oldText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n\n# ---- 1. \"Monsieur le Directeur de la soci\u00e9t\u00e9\u2026\" paragraph ----------------\nReplace-Text \"Monsieur le Directeur de la soci\u00e9t\u00e9\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026\u2026..\" \"Monsieur le Directeur de la soci\u00e9t\u00e9 <societe>\"\n\n# ---- 2. \"Objet\" paragraph -------------------------------------------------\n# \"Objet \" itself stays bold+underlined and untouched; only the remainder of\n# the paragraph (after that label) is replaced.\nReplace-Text \"   : D\u00e9cision d\u2019Institution du permis de recherche n\u00b0 : \u2026\u2026\u2026\u2026. \" \"   : D\u00e9cision d\u2019Institution du permis de recherche n\u00b0 : <Num_PR> \"\n\n# ---- 3. \"Suite \u00e0 votre demande ...\" paragraph -----------------------------\nReplace-Text \"                Suite \u00e0 votre demande d\u2019Institution du permis de recherche n\u00b0 \u2026\u2026\u2026.d\u00e9pos\u00e9e \u00e0 La Direction R\u00e9gionale du D\u00e9partement de l\u2019Energie et des Mines de\u2026.... le \u2026\u2026.., j\u2019ai l\u2019honneur de vous faire parvenir ci-joint la d\u00e9cision cit\u00e9e en objet. \" \"                Suite \u00e0 votre demande d\u2019Institution du permis de recherche n\u00b0 <Num_PR> d\u00e9pos\u00e9e \u00e0 La Direction R\u00e9gionale du D\u00e9partement de l\u2019Energie et des Mines de Marrakech. Le <date>, j\u2019ai l\u2019honneur de vous faire parvenir ci-joint la d\u00e9cision cit\u00e9e en objet. \"\n\n# ---- 4. Closing date/place line -------------------------------------------\nReplace-Text \"\u2026\u2026\u2026\u2026\u2026\u2026, le\u2026\u2026\u2026/\u2026\u2026\u2026/\u2026\u2026\u2026\" \"                                                                   Marrakech, le <date>\"\n\n# Center-align that closing paragraph (was right-aligned).\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.StartsWith(\"                                                                   Marrakech, le\")) {\n        $p.Alignment = 1\n        break\n    }\n}\n"}
